$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '67.190.48'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '3.458.16'
$ws.Range("E3").Value = '  +1.83%  '

$ws.Range("E4").Value = '  -0.13%  '

Set-TextValue "D5" '578.39'
$ws.Range("E5").Value = '  +3.07%  '

Set-TextValue "D6" '187.72'
$ws.Range("E6").Value = '  +6.71%  '

$ws.Range("E7").Value = '  +0.46%  '

$ws.Range("D8").Value = '3.450.14'
$ws.Range("E8").Value = '  +1.89%  '

$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("E10").Value = '  -0.47%  '

Set-TextValue "D11" '0.644'
$ws.Range("E11").Value = '  +0.99%  '

Set-TextValue "D12" '57.97'
$ws.Range("E12").Value = '  +7.67%  '

$ws.Range("E13").Value = '  -1.05%  '

Set-TextValue "D14" '9.49'
$ws.Range("E14").Value = '  +2.91%  '

$ws.Range("D15").Value = '4.001.79'
$ws.Range("E15").Value = '  +1.53%  '

Set-TextValue "D16" '18.96'
$ws.Range("E16").Value = '  +3.11%  '

$ws.Range("D17").Value = '3.454.83'
$ws.Range("E17").Value = '  +1.58%  '

$ws.Range("D18").Value = '67.119.67'
$ws.Range("E18").Value = '  +2.58%  '

$ws.Range("E19").Value = '  -0.44%  '

Set-TextValue "D20" '12.06'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("E21").Value = '  +2.03%  '

Set-TextValue "D22" '486.72'
$ws.Range("E22").Value = '  +5.07%  '

$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D23" '5.50'
$ws.Range("E23").Value = '  +11.45%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D24" '17.18'
$ws.Range("E24").Value = '  +20.55%  '

Set-TextValue "D25" '4.35'
$ws.Range("E25").Value = '  +5.00%  '

Set-TextValue "D26" '89.54'
$ws.Range("E26").Value = '  +2.47%  '

Set-TextValue "D27" '2.97'
$ws.Range("E27").Value = '  +1.23%  '

Set-TextValue "D28" '10.97'
$ws.Range("E28").Value = '  +2.44%  '

Set-TextValue "D29" '9.03'
$ws.Range("E29").Value = '  +3.27%  '

Set-TextValue "D30" '31.27'
$ws.Range("E30").Value = '  +0.51%  '

Set-TextValue "D31" '7.39'
$ws.Range("E31").Value = '  +12.72%  '

Set-TextValue "D32" '602.37'
$ws.Range("E32").Value = '  +4.35%  '

$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D33" '11.84'
$ws.Range("E33").Value = '  +2.81%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D34" '64.90'
$ws.Range("E34").Value = '  +1.85%  '

Set-TextValue "D35" '0.112'
$ws.Range("E35").Value = '  +3.60%  '

$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("E37").Value = '  +2.57%  '

$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0793'
$ws.Range("E38").Value = '  +6.74%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D39" '36.85'
$ws.Range("E39").Value = '  +2.37%  '

Set-TextValue "D40" '0.386'
$ws.Range("E40").Value = '  +3.08%  '

$ws.Range("E41").Value = '  -4.20%  '

$ws.Range("D42").Value = '3.196.76'
$ws.Range("E42").Value = '  +2.94%  '

$ws.Range("E43").Value = '  +2.56%  '

$ws.Range("E44").Value = '  +3.12%  '

Set-TextValue "D45" '2.56'
$ws.Range("E45").Value = '  +4.43%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D46" '0.136'
$ws.Range("E46").Value = '  +1.30%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D47" '3.22'
$ws.Range("E47").Value = '  +1.27%  '

$ws.Range("E48").Value = '  +16.28%  '

Set-TextValue "D49" '0.998'
$ws.Range("E49").Value = '  -0.28%  '

Set-TextValue "D50" '8.68'
$ws.Range("E50").Value = '  +3.03%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D51" '141.45'
$ws.Range("E51").Value = '  +0.80%  '
